$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Location column (C) with the new coordinates for each node row.
$ws.Range("C2").Value = "30.262489, 57.106441"
$ws.Range("C3").Value = "30.292477, 57.089221"
$ws.Range("C4").Value = "30.290032, 57.039864"
$ws.Range("C5").Value = "30.264507, 57.049214"
$ws.Range("C6").Value = "30.307728, 57.098888"
$ws.Range("C7").Value = "30.292568, 57.111453"
$ws.Range("C8").Value = "30.312610, 57.068984"
$ws.Range("C9").Value = "30.311235, 57.033492"
$ws.Range("C10").Value = "30.256178, 57.079154"
$ws.Range("C11").Value = "29.592670, 57.438840"
$ws.Range("C12").Value = "29.100416, 58.369935"
$ws.Range("C13").Value = "28.735959, 57.330483"
$ws.Range("C14").Value = "28.163476, 57.312216"
$ws.Range("C15").Value = "28.812881, 56.547009"
$ws.Range("C16").Value = "29.858242, 56.798987"
$ws.Range("C17").Value = "29.927019, 56.567369"
$ws.Range("C18").Value = "29.430221, 55.672557"
$ws.Range("C19").Value = "30.114571, 55.124759"
$ws.Range("C20").Value = "29.998328, 55.791895"
$ws.Range("C21").Value = "30.810333, 56.586627"
$ws.Range("C22").Value = "30.399309, 56.001607"
$ws.Range("C23").Value = "30.886843, 55.243911"

# Move the active selection to C5.
$ws.Range("C5").Select()

# Touch page setup so paper size / orientation are explicit (A4 portrait).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
